$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = [double]"0.60844"
$ws.Cells.Item(2, 8).Value = [double]"1.82532"
$ws.Cells.Item(2, 9).Value = [double]"0.01418783042133501"
$ws.Cells.Item(2, 10).Value = [double]"0.01418783042133501"
$ws.Cells.Item(2, 13).Value = [double]"1.225147333333333"
$ws.Cells.Item(2, 14).Value = [double]"3.675442"
$ws.Cells.Item(2, 15).Value = [double]"0.2944933560673559"
$ws.Cells.Item(2, 16).Value = [double]"0.2944933560673559"
$ws.Cells.Item(2, 17).Value = [double]"0.7454286434933334"
$ws.Cells.Item(2, 18).Value = [double]"6.708857791440001"
$ws.Cells.Item(2, 19).Value = [double]"0.004178221796093474"
$ws.Cells.Item(2, 20).Value = [double]"0.004178221796093474"
$ws.Cells.Item(3, 7).Value = [double]"0.60844"
$ws.Cells.Item(3, 8).Value = [double]"1.82532"
$ws.Cells.Item(3, 9).Value = [double]"0.01418783042133501"
$ws.Cells.Item(3, 10).Value = [double]"0.01418783042133501"
$ws.Cells.Item(3, 13).Value = [double]"0.005333666666666667"
$ws.Cells.Item(3, 15).Value = [double]"0.001282073881300198"
$ws.Cells.Item(3, 16).Value = [double]"0.001282073881300198"
$ws.Cells.Item(3, 17).Value = [double]"0.003245216146666667"
$ws.Cells.Item(3, 18).Value = [double]"0.02920694532"
$ws.Cells.Item(3, 19).Value = [double]"1.818984681551E-05"
$ws.Cells.Item(3, 20).Value = [double]"1.818984681551E-05"
$ws.Cells.Item(4, 7).Value = [double]"0.60844"
$ws.Cells.Item(4, 8).Value = [double]"1.82532"
$ws.Cells.Item(4, 9).Value = [double]"0.01418783042133501"
$ws.Cells.Item(4, 10).Value = [double]"0.01418783042133501"
$ws.Cells.Item(4, 13).Value = [double]"0.4766303333333334"
$ws.Cells.Item(4, 14).Value = [double]"1.429891"
$ws.Cells.Item(4, 15).Value = [double]"0.1145694584217375"
$ws.Cells.Item(4, 16).Value = [double]"0.1145694584217375"
$ws.Cells.Item(4, 17).Value = [double]"0.2900009600133334"
$ws.Cells.Item(4, 18).Value = [double]"2.610008640120001"
$ws.Cells.Item(4, 19).Value = [double]"0.001625492047551803"
$ws.Cells.Item(4, 20).Value = [double]"0.001625492047551804"
$ws.Cells.Item(5, 7).Value = [double]"0.60844"
$ws.Cells.Item(5, 8).Value = [double]"1.82532"
$ws.Cells.Item(5, 9).Value = [double]"0.01418783042133501"
$ws.Cells.Item(5, 10).Value = [double]"0.01418783042133501"
$ws.Cells.Item(5, 13).Value = [double]"2.453075333333333"
$ws.Cells.Item(5, 14).Value = [double]"7.359226"
$ws.Cells.Item(5, 15).Value = [double]"0.5896551116296064"
$ws.Cells.Item(5, 16).Value = [double]"0.5896551116296064"
$ws.Cells.Item(5, 17).Value = [double]"1.492549155813333"
$ws.Cells.Item(5, 18).Value = [double]"13.43294240232"
$ws.Cells.Item(5, 19).Value = [double]"0.008365926730874218"
$ws.Cells.Item(5, 20).Value = [double]"0.008365926730874218"
$ws.Cells.Item(6, 9).Value = [double]"0.02409258886165303"
$ws.Cells.Item(6, 10).Value = [double]"0.02409258886165303"
$ws.Cells.Item(6, 13).Value = [double]"1.225147333333333"
$ws.Cells.Item(6, 14).Value = [double]"3.675442"
$ws.Cells.Item(6, 15).Value = [double]"0.2944933560673559"
$ws.Cells.Item(6, 16).Value = [double]"0.2944933560673559"
$ws.Cells.Item(6, 17).Value = [double]"1.265824675094667"
$ws.Cells.Item(6, 18).Value = [double]"11.392422075852"
$ws.Cells.Item(6, 19).Value = [double]"0.0070951073502192"
$ws.Cells.Item(6, 20).Value = [double]"0.007095107350219199"
$ws.Cells.Item(7, 9).Value = [double]"0.02409258886165303"
$ws.Cells.Item(7, 10).Value = [double]"0.02409258886165303"
$ws.Cells.Item(7, 13).Value = [double]"0.005333666666666667"
$ws.Cells.Item(7, 15).Value = [double]"0.001282073881300198"
$ws.Cells.Item(7, 16).Value = [double]"0.001282073881300198"
$ws.Cells.Item(7, 17).Value = [double]"0.005510755067333333"
$ws.Cells.Item(7, 18).Value = [double]"0.049596795606"
$ws.Cells.Item(7, 19).Value = [double]"3.088847891242942E-05"
$ws.Cells.Item(7, 20).Value = [double]"3.088847891242942E-05"
$ws.Cells.Item(8, 9).Value = [double]"0.02409258886165303"
$ws.Cells.Item(8, 10).Value = [double]"0.02409258886165303"
$ws.Cells.Item(8, 13).Value = [double]"0.4766303333333334"
$ws.Cells.Item(8, 14).Value = [double]"1.429891"
$ws.Cells.Item(8, 15).Value = [double]"0.1145694584217375"
$ws.Cells.Item(8, 16).Value = [double]"0.1145694584217375"
$ws.Cells.Item(8, 17).Value = [double]"0.4924554136606668"
$ws.Cells.Item(8, 18).Value = [double]"4.432098722946"
$ws.Cells.Item(8, 19).Value = [double]"0.002760274857857173"
$ws.Cells.Item(8, 20).Value = [double]"0.002760274857857173"
$ws.Cells.Item(9, 9).Value = [double]"0.02409258886165303"
$ws.Cells.Item(9, 10).Value = [double]"0.02409258886165303"
$ws.Cells.Item(9, 13).Value = [double]"2.453075333333333"
$ws.Cells.Item(9, 14).Value = [double]"7.359226"
$ws.Cells.Item(9, 15).Value = [double]"0.5896551116296064"
$ws.Cells.Item(9, 16).Value = [double]"0.5896551116296064"
$ws.Cells.Item(9, 17).Value = [double]"2.534522340550667"
$ws.Cells.Item(9, 18).Value = [double]"22.810701064956"
$ws.Cells.Item(9, 19).Value = [double]"0.01420631817466423"
$ws.Cells.Item(9, 20).Value = [double]"0.01420631817466423"
$ws.Cells.Item(10, 7).Value = [double]"6.042074333333333"
$ws.Cells.Item(10, 8).Value = [double]"18.126223"
$ws.Cells.Item(10, 9).Value = [double]"0.1408913385616233"
$ws.Cells.Item(10, 10).Value = [double]"0.1408913385616233"
$ws.Cells.Item(10, 13).Value = [double]"1.225147333333333"
$ws.Cells.Item(10, 14).Value = [double]"3.675442"
$ws.Cells.Item(10, 15).Value = [double]"0.2944933560673559"
$ws.Cells.Item(10, 16).Value = [double]"0.2944933560673559"
$ws.Cells.Item(10, 17).Value = [double]"7.402431257285111"
$ws.Cells.Item(10, 18).Value = [double]"66.621881315566"
$ws.Cells.Item(10, 19).Value = [double]"0.04149156313383453"
$ws.Cells.Item(10, 20).Value = [double]"0.04149156313383453"
$ws.Cells.Item(11, 7).Value = [double]"6.042074333333333"
$ws.Cells.Item(11, 8).Value = [double]"18.126223"
$ws.Cells.Item(11, 9).Value = [double]"0.1408913385616233"
$ws.Cells.Item(11, 10).Value = [double]"0.1408913385616233"
$ws.Cells.Item(11, 13).Value = [double]"0.005333666666666667"
$ws.Cells.Item(11, 15).Value = [double]"0.001282073881300198"
$ws.Cells.Item(11, 16).Value = [double]"0.001282073881300198"
$ws.Cells.Item(11, 17).Value = [double]"0.03222641046922222"
$ws.Cells.Item(11, 18).Value = [double]"0.290037694223"
$ws.Cells.Item(11, 19).Value = [double]"0.0001806331052712807"
$ws.Cells.Item(11, 20).Value = [double]"0.0001806331052712807"
$ws.Cells.Item(12, 7).Value = [double]"6.042074333333333"
$ws.Cells.Item(12, 8).Value = [double]"18.126223"
$ws.Cells.Item(12, 9).Value = [double]"0.1408913385616233"
$ws.Cells.Item(12, 10).Value = [double]"0.1408913385616233"
$ws.Cells.Item(12, 13).Value = [double]"0.4766303333333334"
$ws.Cells.Item(12, 14).Value = [double]"1.429891"
$ws.Cells.Item(12, 15).Value = [double]"0.1145694584217375"
$ws.Cells.Item(12, 16).Value = [double]"0.1145694584217375"
$ws.Cells.Item(12, 17).Value = [double]"2.879835903521445"
$ws.Cells.Item(12, 18).Value = [double]"25.918523131693"
$ws.Cells.Item(12, 19).Value = [double]"0.01614184435531884"
$ws.Cells.Item(12, 20).Value = [double]"0.01614184435531884"
$ws.Cells.Item(13, 7).Value = [double]"6.042074333333333"
$ws.Cells.Item(13, 8).Value = [double]"18.126223"
$ws.Cells.Item(13, 9).Value = [double]"0.1408913385616233"
$ws.Cells.Item(13, 10).Value = [double]"0.1408913385616233"
$ws.Cells.Item(13, 13).Value = [double]"2.453075333333333"
$ws.Cells.Item(13, 14).Value = [double]"7.359226"
$ws.Cells.Item(13, 15).Value = [double]"0.5896551116296064"
$ws.Cells.Item(13, 16).Value = [double]"0.5896551116296064"
$ws.Cells.Item(13, 17).Value = [double]"14.82166350926645"
$ws.Cells.Item(13, 18).Value = [double]"133.394971583398"
$ws.Cells.Item(13, 19).Value = [double]"0.08307729796719866"
$ws.Cells.Item(13, 20).Value = [double]"0.08307729796719866"
$ws.Cells.Item(14, 7).Value = [double]"0.4072233333333333"
$ws.Cells.Item(14, 8).Value = [double]"1.22167"
$ws.Cells.Item(14, 9).Value = [double]"0.009495785281940885"
$ws.Cells.Item(14, 10).Value = [double]"0.009495785281940885"
$ws.Cells.Item(14, 13).Value = [double]"1.225147333333333"
$ws.Cells.Item(14, 14).Value = [double]"3.675442"
$ws.Cells.Item(14, 15).Value = [double]"0.2944933560673559"
$ws.Cells.Item(14, 16).Value = [double]"0.2944933560673559"
$ws.Cells.Item(14, 17).Value = [double]"0.4989085809044445"
$ws.Cells.Item(14, 18).Value = [double]"4.49017722814"
$ws.Cells.Item(14, 19).Value = [double]"0.002796445676173775"
$ws.Cells.Item(14, 20).Value = [double]"0.002796445676173775"
$ws.Cells.Item(15, 7).Value = [double]"0.4072233333333333"
$ws.Cells.Item(15, 8).Value = [double]"1.22167"
$ws.Cells.Item(15, 9).Value = [double]"0.009495785281940885"
$ws.Cells.Item(15, 10).Value = [double]"0.009495785281940885"
$ws.Cells.Item(15, 13).Value = [double]"0.005333666666666667"
$ws.Cells.Item(15, 15).Value = [double]"0.001282073881300198"
$ws.Cells.Item(15, 16).Value = [double]"0.001282073881300198"
$ws.Cells.Item(15, 17).Value = [double]"0.002171993518888889"
$ws.Cells.Item(15, 18).Value = [double]"0.01954794167"
$ws.Cells.Item(15, 19).Value = [double]"1.217429829241125E-05"
$ws.Cells.Item(15, 20).Value = [double]"1.217429829241125E-05"
$ws.Cells.Item(16, 7).Value = [double]"0.4072233333333333"
$ws.Cells.Item(16, 8).Value = [double]"1.22167"
$ws.Cells.Item(16, 9).Value = [double]"0.009495785281940885"
$ws.Cells.Item(16, 10).Value = [double]"0.009495785281940885"
$ws.Cells.Item(16, 13).Value = [double]"0.4766303333333334"
$ws.Cells.Item(16, 14).Value = [double]"1.429891"
$ws.Cells.Item(16, 15).Value = [double]"0.1145694584217375"
$ws.Cells.Item(16, 16).Value = [double]"0.1145694584217375"
$ws.Cells.Item(16, 17).Value = [double]"0.1940949931077778"
$ws.Cells.Item(16, 18).Value = [double]"1.74685493797"
$ws.Cells.Item(16, 19).Value = [double]"0.001087926977041073"
$ws.Cells.Item(16, 20).Value = [double]"0.001087926977041073"
$ws.Cells.Item(17, 7).Value = [double]"0.4072233333333333"
$ws.Cells.Item(17, 8).Value = [double]"1.22167"
$ws.Cells.Item(17, 9).Value = [double]"0.009495785281940885"
$ws.Cells.Item(17, 10).Value = [double]"0.009495785281940885"
$ws.Cells.Item(17, 13).Value = [double]"2.453075333333333"
$ws.Cells.Item(17, 14).Value = [double]"7.359226"
$ws.Cells.Item(17, 15).Value = [double]"0.5896551116296064"
$ws.Cells.Item(17, 16).Value = [double]"0.5896551116296064"
$ws.Cells.Item(17, 17).Value = [double]"0.9989495141577778"
$ws.Cells.Item(17, 18).Value = [double]"8.990545627420001"
$ws.Cells.Item(17, 19).Value = [double]"0.005599238330433626"
$ws.Cells.Item(17, 20).Value = [double]"0.005599238330433626"
$ws.Cells.Item(18, 7).Value = [double]"31.00247266666667"
$ws.Cells.Item(18, 8).Value = [double]"93.007418"
$ws.Cells.Item(18, 9).Value = [double]"0.72292719879814"
$ws.Cells.Item(18, 10).Value = [double]"0.72292719879814"
$ws.Cells.Item(18, 13).Value = [double]"1.225147333333333"
$ws.Cells.Item(18, 14).Value = [double]"3.675442"
$ws.Cells.Item(18, 15).Value = [double]"0.2944933560673559"
$ws.Cells.Item(18, 16).Value = [double]"0.2944933560673559"
$ws.Cells.Item(18, 17).Value = [double]"37.98259671430623"
$ws.Cells.Item(18, 18).Value = [double]"341.843370428756"
$ws.Cells.Item(18, 19).Value = [double]"0.2128972569664369"
$ws.Cells.Item(18, 20).Value = [double]"0.2128972569664369"
$ws.Cells.Item(19, 7).Value = [double]"31.00247266666667"
$ws.Cells.Item(19, 8).Value = [double]"93.007418"
$ws.Cells.Item(19, 9).Value = [double]"0.72292719879814"
$ws.Cells.Item(19, 10).Value = [double]"0.72292719879814"
$ws.Cells.Item(19, 13).Value = [double]"0.005333666666666667"
$ws.Cells.Item(19, 15).Value = [double]"0.001282073881300198"
$ws.Cells.Item(19, 16).Value = [double]"0.001282073881300198"
$ws.Cells.Item(19, 17).Value = [double]"0.1653568550464445"
$ws.Cells.Item(19, 18).Value = [double]"1.488211695418"
$ws.Cells.Item(19, 19).Value = [double]"0.0009268460796606113"
$ws.Cells.Item(19, 20).Value = [double]"0.0009268460796606113"
$ws.Cells.Item(20, 7).Value = [double]"31.00247266666667"
$ws.Cells.Item(20, 8).Value = [double]"93.007418"
$ws.Cells.Item(20, 9).Value = [double]"0.72292719879814"
$ws.Cells.Item(20, 10).Value = [double]"0.72292719879814"
$ws.Cells.Item(20, 13).Value = [double]"0.4766303333333334"
$ws.Cells.Item(20, 14).Value = [double]"1.429891"
$ws.Cells.Item(20, 15).Value = [double]"0.1145694584217375"
$ws.Cells.Item(20, 16).Value = [double]"0.1145694584217375"
$ws.Cells.Item(20, 17).Value = [double]"14.77671888127089"
$ws.Cells.Item(20, 18).Value = [double]"132.990469931438"
$ws.Cells.Item(20, 19).Value = [double]"0.08282537764464666"
$ws.Cells.Item(20, 20).Value = [double]"0.08282537764464667"
$ws.Cells.Item(21, 7).Value = [double]"31.00247266666667"
$ws.Cells.Item(21, 8).Value = [double]"93.007418"
$ws.Cells.Item(21, 9).Value = [double]"0.72292719879814"
$ws.Cells.Item(21, 10).Value = [double]"0.72292719879814"
$ws.Cells.Item(21, 13).Value = [double]"2.453075333333333"
$ws.Cells.Item(21, 14).Value = [double]"7.359226"
$ws.Cells.Item(21, 15).Value = [double]"0.5896551116296064"
$ws.Cells.Item(21, 16).Value = [double]"0.5896551116296064"
$ws.Cells.Item(21, 17).Value = [double]"76.05140097094089"
$ws.Cells.Item(21, 18).Value = [double]"684.462608738468"
$ws.Cells.Item(21, 19).Value = [double]"0.4262777181073959"
$ws.Cells.Item(21, 20).Value = [double]"0.4262777181073959"
$ws.Cells.Item(22, 7).Value = [double]"3.791227666666666"
$ws.Cells.Item(22, 8).Value = [double]"11.373683"
$ws.Cells.Item(22, 9).Value = [double]"0.08840525807530777"
$ws.Cells.Item(22, 10).Value = [double]"0.08840525807530777"
$ws.Cells.Item(22, 13).Value = [double]"1.225147333333333"
$ws.Cells.Item(22, 14).Value = [double]"3.675442"
$ws.Cells.Item(22, 15).Value = [double]"0.2944933560673559"
$ws.Cells.Item(22, 16).Value = [double]"0.2944933560673559"
$ws.Cells.Item(22, 17).Value = [double]"4.644812465876222"
$ws.Cells.Item(22, 18).Value = [double]"41.803312192886"
$ws.Cells.Item(22, 19).Value = [double]"0.02603476114459811"
$ws.Cells.Item(22, 20).Value = [double]"0.02603476114459811"
$ws.Cells.Item(23, 7).Value = [double]"3.791227666666666"
$ws.Cells.Item(23, 8).Value = [double]"11.373683"
$ws.Cells.Item(23, 9).Value = [double]"0.08840525807530777"
$ws.Cells.Item(23, 10).Value = [double]"0.08840525807530777"
$ws.Cells.Item(23, 13).Value = [double]"0.005333666666666667"
$ws.Cells.Item(23, 15).Value = [double]"0.001282073881300198"
$ws.Cells.Item(23, 16).Value = [double]"0.001282073881300198"
$ws.Cells.Item(23, 17).Value = [double]"0.02022114463144445"
$ws.Cells.Item(23, 18).Value = [double]"0.181990301683"
$ws.Cells.Item(23, 19).Value = [double]"0.0001133420723479555"
$ws.Cells.Item(23, 20).Value = [double]"0.0001133420723479555"
$ws.Cells.Item(24, 7).Value = [double]"3.791227666666666"
$ws.Cells.Item(24, 8).Value = [double]"11.373683"
$ws.Cells.Item(24, 9).Value = [double]"0.08840525807530777"
$ws.Cells.Item(24, 10).Value = [double]"0.08840525807530777"
$ws.Cells.Item(24, 13).Value = [double]"0.4766303333333334"
$ws.Cells.Item(24, 14).Value = [double]"1.429891"
$ws.Cells.Item(24, 15).Value = [double]"0.1145694584217375"
$ws.Cells.Item(24, 16).Value = [double]"0.1145694584217375"
$ws.Cells.Item(24, 17).Value = [double]"1.807014106505889"
$ws.Cells.Item(24, 18).Value = [double]"16.263126958553"
$ws.Cells.Item(24, 19).Value = [double]"0.01012854253932195"
$ws.Cells.Item(24, 20).Value = [double]"0.01012854253932195"
$ws.Cells.Item(25, 7).Value = [double]"3.791227666666666"
$ws.Cells.Item(25, 8).Value = [double]"11.373683"
$ws.Cells.Item(25, 9).Value = [double]"0.08840525807530777"
$ws.Cells.Item(25, 10).Value = [double]"0.08840525807530777"
$ws.Cells.Item(25, 13).Value = [double]"2.453075333333333"
$ws.Cells.Item(25, 14).Value = [double]"7.359226"
$ws.Cells.Item(25, 15).Value = [double]"0.5896551116296064"
$ws.Cells.Item(25, 16).Value = [double]"0.5896551116296064"
$ws.Cells.Item(25, 17).Value = [double]"9.30016707215089"
$ws.Cells.Item(25, 18).Value = [double]"83.701503649358"
$ws.Cells.Item(25, 19).Value = [double]"0.05212861231903976"
$ws.Cells.Item(25, 20).Value = [double]"0.05212861231903976"
